# Appends a new "How Moolah Capital Can Help Beginners" section (plus a
# restated Disclaimer paragraph) to the very end of the document, after
# the last (empty) paragraph and before the section break.
#
# Non-ASCII characters (bullet "•", em dash "—") are built with [char] /
# codepoints rather than typed literally, since literal multi-byte UTF-8
# in the script source does not always survive the COM-interop bridge.

$d = $word.ActiveDocument

$bullet = [char]0x2022   # "•"
$dash   = [char]0x2014   # "—"

$newParagraphs = @(
    "",
    "How Moolah Capital Can Help Beginners",
    "${bullet}Market Index Fund $dash simple diversified exposure as you learn",
    "${bullet}AlphaGlobal Momentum Fund $dash rules-based trend exposure when you're ready to scale",
    "${bullet}AlphaGlobal Yield Fund $dash income via staking/lending with risk controls",
    "${bullet}GenAI Funds $dash build or mirror AI-driven strategies as you get comfortable",
    "",
    "Disclaimer: This guide is UK-oriented; tax rules vary by country${dash}seek professional advice. Nothing here is financial advice."
)

# Anchor on the last paragraph currently in the document (the trailing
# empty paragraph right before the sectPr) and insert each new paragraph
# after it in order, so formatting (font/size) is inherited from that
# paragraph's run, matching the rest of the document.
$insertAfter = $d.Paragraphs.Last.Range

foreach ($text in $newParagraphs) {
    $insertAfter.InsertParagraphAfter()
    $newPara = $d.Paragraphs.Last
    if ($text -ne "") {
        $newPara.Range.Text = $text
    }
    $insertAfter = $newPara.Range
}

Write-Host "Paragraphs now:" $d.Paragraphs.Count
